$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BP-813: Affiliate Mapping for True Independent Stations
# Apply the same bordered/body formatting used by the rest of the table
# to columns I and J (previously unbordered).
$ws.Range("B1:B5").Copy()
$ws.Range("I1:J5").PasteSpecial(-4122)

# Rename "Affiliation Mismatch Note" header -> "IsTrueIND"
# Rename "SalesGroupName" header -> "RepFirm"
$ws.Range("H1").Value = "IsTrueIND"
$ws.Range("J1").Value = "RepFirm"

$ws.Range("A1").Select()
